# Update column F (dSF) values for specific rows per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    4  = -5
    11 = 1
    14 = -12
    15 = 1
    17 = 1
    25 = 3
    28 = 4
    32 = -2
    34 = -4
    36 = 9
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
